$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.90715
$ws.Cells.Item(2, 8).Value = 2.72145
$ws.Cells.Item(2, 9).Value = 0.01717809939998381
$ws.Cells.Item(2, 10).Value = 0.01717809939998381
$ws.Cells.Item(2, 13).Value = 3.443466666666667
$ws.Cells.Item(2, 14).Value = 10.3304
$ws.Cells.Item(2, 15).Value = 0.07502986933839939
$ws.Cells.Item(2, 16).Value = 0.07502986933839939
$ws.Cells.Item(2, 17).Value = 3.123740786666667
$ws.Cells.Item(2, 18).Value = 28.11366708
$ws.Cells.Item(2, 19).Value = 0.001288870553462822
$ws.Cells.Item(2, 20).Value = 0.001288870553462822
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.90715
$ws.Cells.Item(3, 8).Value = 2.72145
$ws.Cells.Item(3, 9).Value = 0.01717809939998381
$ws.Cells.Item(3, 10).Value = 0.01717809939998381
$ws.Cells.Item(3, 15).Value = 0.08813227911805739
$ws.Cells.Item(3, 16).Value = 0.08813227911805739
$ws.Cells.Item(3, 17).Value = 3.669237296166667
$ws.Cells.Item(3, 18).Value = 33.0231356655
$ws.Cells.Item(3, 19).Value = 0.001513945051037108
$ws.Cells.Item(3, 20).Value = 0.001513945051037107
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.90715
$ws.Cells.Item(4, 8).Value = 2.72145
$ws.Cells.Item(4, 9).Value = 0.01717809939998381
$ws.Cells.Item(4, 10).Value = 0.01717809939998381
$ws.Cells.Item(4, 13).Value = 2.521553333333333
$ws.Cells.Item(4, 14).Value = 7.56466
$ws.Cells.Item(4, 15).Value = 0.054942253096629
$ws.Cells.Item(4, 16).Value = 0.054942253096629
$ws.Cells.Item(4, 17).Value = 2.287427106333333
$ws.Cells.Item(4, 18).Value = 20.586843957
$ws.Cells.Item(4, 19).Value = 0.0009438034849529614
$ws.Cells.Item(4, 20).Value = 0.0009438034849529611
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.90715
$ws.Cells.Item(5, 8).Value = 2.72145
$ws.Cells.Item(5, 9).Value = 0.01717809939998381
$ws.Cells.Item(5, 10).Value = 0.01717809939998381
$ws.Cells.Item(5, 13).Value = 35.88479433333333
$ws.Cells.Item(5, 14).Value = 107.654383
$ws.Cells.Item(5, 15).Value = 0.7818955984469143
$ws.Cells.Item(5, 16).Value = 0.7818955984469143
$ws.Cells.Item(5, 17).Value = 32.55289117948333
$ws.Cells.Item(5, 18).Value = 292.97602061535
$ws.Cells.Item(5, 19).Value = 0.01343148031053092
$ws.Cells.Item(5, 20).Value = 0.01343148031053092
$ws.Cells.Item(6, 9).Value = 0.0237180037344858
$ws.Cells.Item(6, 10).Value = 0.0237180037344858
$ws.Cells.Item(6, 13).Value = 3.443466666666667
$ws.Cells.Item(6, 14).Value = 10.3304
$ws.Cells.Item(6, 15).Value = 0.07502986933839939
$ws.Cells.Item(6, 16).Value = 0.07502986933839939
$ws.Cells.Item(6, 17).Value = 4.312985617244445
$ws.Cells.Item(6, 18).Value = 38.8168705552
$ws.Cells.Item(6, 19).Value = 0.001779558721166138
$ws.Cells.Item(6, 20).Value = 0.001779558721166138
$ws.Cells.Item(7, 9).Value = 0.0237180037344858
$ws.Cells.Item(7, 10).Value = 0.0237180037344858
$ws.Cells.Item(7, 15).Value = 0.08813227911805739
$ws.Cells.Item(7, 16).Value = 0.08813227911805739
$ws.Cells.Item(7, 19).Value = 0.00209032172525083
$ws.Cells.Item(7, 20).Value = 0.00209032172525083
$ws.Cells.Item(8, 9).Value = 0.0237180037344858
$ws.Cells.Item(8, 10).Value = 0.0237180037344858
$ws.Cells.Item(8, 13).Value = 2.521553333333333
$ws.Cells.Item(8, 14).Value = 7.56466
$ws.Cells.Item(8, 15).Value = 0.054942253096629
$ws.Cells.Item(8, 16).Value = 0.054942253096629
$ws.Cells.Item(8, 17).Value = 3.158277489675556
$ws.Cells.Item(8, 18).Value = 28.42449740708
$ws.Cells.Item(8, 19).Value = 0.001303120564126911
$ws.Cells.Item(8, 20).Value = 0.00130312056412691
$ws.Cells.Item(9, 9).Value = 0.0237180037344858
$ws.Cells.Item(9, 10).Value = 0.0237180037344858
$ws.Cells.Item(9, 13).Value = 35.88479433333333
$ws.Cells.Item(9, 14).Value = 107.654383
$ws.Cells.Item(9, 15).Value = 0.7818955984469143
$ws.Cells.Item(9, 16).Value = 0.7818955984469143
$ws.Cells.Item(9, 17).Value = 44.94615944322823
$ws.Cells.Item(9, 18).Value = 404.515434989054
$ws.Cells.Item(9, 19).Value = 0.01854500272394192
$ws.Cells.Item(9, 20).Value = 0.01854500272394192
$ws.Cells.Item(10, 7).Value = 2.247832333333333
$ws.Cells.Item(10, 8).Value = 6.743497
$ws.Cells.Item(10, 9).Value = 0.04256571378106988
$ws.Cells.Item(10, 10).Value = 0.04256571378106987
$ws.Cells.Item(10, 13).Value = 3.443466666666667
$ws.Cells.Item(10, 14).Value = 10.3304
$ws.Cells.Item(10, 15).Value = 0.07502986933839939
$ws.Cells.Item(10, 16).Value = 0.07502986933839939
$ws.Cells.Item(10, 17).Value = 7.740335712088889
$ws.Cells.Item(10, 18).Value = 69.66302140880001
$ws.Cells.Item(10, 19).Value = 0.003193699943289379
$ws.Cells.Item(10, 20).Value = 0.003193699943289379
$ws.Cells.Item(11, 7).Value = 2.247832333333333
$ws.Cells.Item(11, 8).Value = 6.743497
$ws.Cells.Item(11, 9).Value = 0.04256571378106988
$ws.Cells.Item(11, 10).Value = 0.04256571378106987
$ws.Cells.Item(11, 15).Value = 0.08813227911805739
$ws.Cells.Item(11, 16).Value = 0.08813227911805739
$ws.Cells.Item(11, 17).Value = 9.092024729092222
$ws.Cells.Item(11, 18).Value = 81.82822256182999
$ws.Cells.Item(11, 19).Value = 0.003751413367812593
$ws.Cells.Item(11, 20).Value = 0.003751413367812592
$ws.Cells.Item(12, 7).Value = 2.247832333333333
$ws.Cells.Item(12, 8).Value = 6.743497
$ws.Cells.Item(12, 9).Value = 0.04256571378106988
$ws.Cells.Item(12, 10).Value = 0.04256571378106987
$ws.Cells.Item(12, 13).Value = 2.521553333333333
$ws.Cells.Item(12, 14).Value = 7.56466
$ws.Cells.Item(12, 15).Value = 0.054942253096629
$ws.Cells.Item(12, 16).Value = 0.054942253096629
$ws.Cells.Item(12, 17).Value = 5.66802911289111
$ws.Cells.Item(12, 18).Value = 51.01226201602
$ws.Cells.Item(12, 19).Value = 0.00233865621979821
$ws.Cells.Item(12, 20).Value = 0.00233865621979821
$ws.Cells.Item(13, 7).Value = 2.247832333333333
$ws.Cells.Item(13, 8).Value = 6.743497
$ws.Cells.Item(13, 9).Value = 0.04256571378106988
$ws.Cells.Item(13, 10).Value = 0.04256571378106987
$ws.Cells.Item(13, 13).Value = 35.88479433333333
$ws.Cells.Item(13, 14).Value = 107.654383
$ws.Cells.Item(13, 15).Value = 0.7818955984469143
$ws.Cells.Item(13, 16).Value = 0.7818955984469143
$ws.Cells.Item(13, 17).Value = 80.66300097748343
$ws.Cells.Item(13, 18).Value = 725.9670087973509
$ws.Cells.Item(13, 19).Value = 0.0332819442501697
$ws.Cells.Item(13, 20).Value = 0.03328194425016969
$ws.Cells.Item(14, 7).Value = 48.40102466666667
$ws.Cells.Item(14, 8).Value = 145.203074
$ws.Cells.Item(14, 9).Value = 0.9165381830844606
$ws.Cells.Item(14, 10).Value = 0.9165381830844604
$ws.Cells.Item(14, 13).Value = 3.443466666666667
$ws.Cells.Item(14, 14).Value = 10.3304
$ws.Cells.Item(14, 15).Value = 0.07502986933839939
$ws.Cells.Item(14, 16).Value = 0.07502986933839939
$ws.Cells.Item(14, 17).Value = 166.6673150721778
$ws.Cells.Item(14, 18).Value = 1500.0058356496
$ws.Cells.Item(14, 19).Value = 0.06876774012048105
$ws.Cells.Item(14, 20).Value = 0.06876774012048105
$ws.Cells.Item(15, 7).Value = 48.40102466666667
$ws.Cells.Item(15, 8).Value = 145.203074
$ws.Cells.Item(15, 9).Value = 0.9165381830844606
$ws.Cells.Item(15, 10).Value = 0.9165381830844604
$ws.Cells.Item(15, 15).Value = 0.08813227911805739
$ws.Cells.Item(15, 16).Value = 0.08813227911805739
$ws.Cells.Item(15, 17).Value = 195.7723032349845
$ws.Cells.Item(15, 18).Value = 1761.95072911486
$ws.Cells.Item(15, 19).Value = 0.08077659897395686
$ws.Cells.Item(15, 20).Value = 0.08077659897395685
$ws.Cells.Item(16, 7).Value = 48.40102466666667
$ws.Cells.Item(16, 8).Value = 145.203074
$ws.Cells.Item(16, 9).Value = 0.9165381830844606
$ws.Cells.Item(16, 10).Value = 0.9165381830844604
$ws.Cells.Item(16, 13).Value = 2.521553333333333
$ws.Cells.Item(16, 14).Value = 7.56466
$ws.Cells.Item(16, 15).Value = 0.054942253096629
$ws.Cells.Item(16, 16).Value = 0.054942253096629
$ws.Cells.Item(16, 17).Value = 122.0457650849822
$ws.Cells.Item(16, 18).Value = 1098.41188576484
$ws.Cells.Item(16, 19).Value = 0.05035667282775092
$ws.Cells.Item(16, 20).Value = 0.05035667282775091
$ws.Cells.Item(17, 7).Value = 48.40102466666667
$ws.Cells.Item(17, 8).Value = 145.203074
$ws.Cells.Item(17, 9).Value = 0.9165381830844606
$ws.Cells.Item(17, 10).Value = 0.9165381830844604
$ws.Cells.Item(17, 13).Value = 35.88479433333333
$ws.Cells.Item(17, 14).Value = 107.654383
$ws.Cells.Item(17, 15).Value = 0.7818955984469143
$ws.Cells.Item(17, 16).Value = 0.7818955984469143
$ws.Cells.Item(17, 17).Value = 1736.860815685927
$ws.Cells.Item(17, 18).Value = 15631.74734117334
$ws.Cells.Item(17, 19).Value = 0.7166371711622718
$ws.Cells.Item(17, 20).Value = 0.7166371711622717
